$d = $word.ActiveDocument

# The document contains three placeholder "<id>p096v_aN</id>" markers (one per
# <div> entry) that were split across three separate runs each:
#   run1: "<id>"        (Courier New, color 7f6000, sz 18)
#   run2: "p096v_aN"     (default font, color 000000)
#   run3: "</id>"       (Courier New, color 7f6000, sz 18)
#
# Replace each with the final, newly-downloaded id text collapsed into a
# single run (formatting of the surrounding <id>/</id> runs is preserved
# automatically since Find/Replace here only touches the matched text).

$d.Content.Find.Execute("<id>p096v_a1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p096v_1</id>", 2)
$d.Content.Find.Execute("<id>p096v_a2</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p096v_2</id>", 2)
$d.Content.Find.Execute("<id>p096v_a3</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p096v_3</id>", 2)
